$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Directorio")

# --- Row 21: "4." / "Entregas y trámite" -------------------------------
# C21 is a genuine number, safe to set directly.
$ws.Range("C21").Value = 1
# D21 is plain text (not numeric-looking), safe to set directly.
$ws.Range("D21").Value = "Entregas y trámite"

# A21 = "4." looks numeric to Excel's auto-detection, so write it as a
# formula returning the literal string, then collapse it down to a plain
# value via copy / paste-special-values. This keeps the original cell
# style (s="1") instead of forcing a new "text" number format style.
$ws.Range("A21").Formula = "=""4."""

# --- Row 22: "4.1." / Zoho link row -------------------------------------
$ws.Range("B22").Formula = "=""4."""
$ws.Range("A22").Formula = "=""4.1."""

$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "Enlace Zoho - Avance de entregas"
$ws.Range("E22").Value = "url "
$ws.Range("F22").Value = 'Enlace Zoho permalink  "Avance de entregas"'
$ws.Range("G22").Value = "https://analytics.zoho.com/open-view/2403793000018042241/72f498ac7974ce7ef70f46e19c1b4a33"

# Collapse the helper formulas in A21/A22/B22 down to literal text values
# while preserving their original cell formatting.
$rng = $ws.Range("A21:H22")
$rng.Copy()
$rng.PasteSpecial(-4163)

# Add the hyperlink for G22, then restore its formatting to match the
# other URL cells in the table (style used by G5/G7/G9/G15/G18/G19/G20).
$ws.Hyperlinks.Add($ws.Range("G22"), "https://analytics.zoho.com/open-view/2403793000018042241/72f498ac7974ce7ef70f46e19c1b4a33")
$ws.Range("G5").Copy()
$ws.Range("G22").PasteSpecial(-4122)

# --- Selection / view state ---------------------------------------------
$ws.Range("D23").Select()
